$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Title cell's value from "CURA Healthcare" to "CURA Healthcare Service"
$ws.Range("G4").Value = "CURA Healthcare Service"

# Update the active selection to match the saved view state
$ws.Range("F14").Select()
